# "Generate Report for Handoff"
#
# For the files that were previously pending handoff (Priority = "low"),
# a new handoff xliff run has just completed:
#   - Priority is bumped from "low" to "ht" (high-temp / handed-off)
#   - The "Latest Handoff Datetime" is refreshed to the new generation time
#
# This touches rows 4-7 (the four files whose Priority was "low") on both
# the "zh-cn" and "de-de" language sheets, plus the rollup "Latest HO Xliff
# Generate Date" column on the Overview sheet (which mirrors the de-de
# handoff time for these rows).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-11-02 05:45:37"
$zhcn.Range("H5").Value = "2016-11-02 05:45:37"
$zhcn.Range("H6").Value = "2016-11-02 05:45:37"
$zhcn.Range("H7").Value = "2016-11-02 05:45:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4").Value = "2016-11-02 05:45:52"
$dede.Range("H5").Value = "2016-11-02 05:45:52"
$dede.Range("H6").Value = "2016-11-02 05:45:52"
$dede.Range("H7").Value = "2016-11-02 05:45:52"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-11-02 05:45:52"
$overview.Range("G5").Value = "2016-11-02 05:45:52"
$overview.Range("G6").Value = "2016-11-02 05:45:52"
$overview.Range("G7").Value = "2016-11-02 05:45:52"
